# Adds the 2025-09-04 daily update to the violent-crime-full-year workbook.
# For each affected sheet, update the 2025 (column L) running totals, and the
# couple of prior-year (J/K) corrections that shipped in the same data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range('L2').Value = 44
$ws.Range('L7').Value = 124

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('L2').Value = 168
$ws.Range('L4').Value = 35
$ws.Range('L7').Value = 499

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('L3').Value = 339
$ws.Range('L4').Value = 75
$ws.Range('L6').Value = 270
$ws.Range('L7').Value = 1011

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range('L2').Value = 90
$ws.Range('L3').Value = 76
$ws.Range('L4').Value = 19
$ws.Range('L7').Value = 244

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range('L2').Value = 34
$ws.Range('L7').Value = 75

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('L2').Value = 124
$ws.Range('L7').Value = 499
$ws.Range('L8').Value = 1011
$ws.Range('L11').Value = 244
$ws.Range('L14').Value = 75
$ws.Range('L18').Value = 104
$ws.Range('L19').Value = 417
$ws.Range('L20').Value = 386
$ws.Range('L27').Value = 133
$ws.Range('L29').Value = 840
$ws.Range('L33').Value = 687
$ws.Range('L34').Value = 90
$ws.Range('L36').Value = 196
$ws.Range('L38').Value = 18
$ws.Range('L42').Value = 493
$ws.Range('L43').Value = 109
$ws.Range('L51').Value = 184
$ws.Range('L53').Value = 174
$ws.Range('L54').Value = 316
$ws.Range('L61').Value = 18
$ws.Range('J63').Value = 224
$ws.Range('L63').Value = 44
$ws.Range('L64').Value = 103
$ws.Range('L65').Value = 290
$ws.Range('L67').Value = 520
$ws.Range('K73').Value = 242
$ws.Range('L73').Value = 121
$ws.Range('L76').Value = 236
$ws.Range('L77').Value = 103
$ws.Range('L78').Value = 199
$ws.Range('L79').Value = 400
$ws.Range('L84').Value = 146
$ws.Range('L85').Value = 777
$ws.Range('L89').Value = 215
$ws.Range('L90').Value = 150
$ws.Range('L94').Value = 190
$ws.Range('L95').Value = 203
$ws.Range('L96').Value = 171
$ws.Range('L99').Value = 259
$ws.Range('J101').Value = 29349
$ws.Range('K101').Value = 27571
$ws.Range('L101').Value = 15137

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range('L2').Value = 38
$ws.Range('L3').Value = 37
$ws.Range('L6').Value = 16
$ws.Range('L7').Value = 104

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('L2').Value = 147
$ws.Range('L3').Value = 128
$ws.Range('L6').Value = 121
$ws.Range('L7').Value = 417

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range('L3').Value = 126
$ws.Range('L6').Value = 105
$ws.Range('L7').Value = 386

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('L2').Value = 4569
$ws.Range('L3').Value = 4895
$ws.Range('J4').Value = 1873
$ws.Range('K4').Value = 1779
$ws.Range('L4').Value = 1209
$ws.Range('L5').Value = 283
$ws.Range('L6').Value = 4181
$ws.Range('J7').Value = 29349
$ws.Range('K7').Value = 27571
$ws.Range('L7').Value = 15137

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range('L2').Value = 35
$ws.Range('L7').Value = 133

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('L2').Value = 252
$ws.Range('L3').Value = 316
$ws.Range('L6').Value = 220
$ws.Range('L7').Value = 840

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('L2').Value = 188
$ws.Range('L5').Value = 14
$ws.Range('L7').Value = 687

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range('L6').Value = 29
$ws.Range('L7').Value = 90

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range('L2').Value = 74
$ws.Range('L7').Value = 196

$ws = $wb.Worksheets.Item('Grant Park')
$ws.Range('L5').Value = 7
$ws.Range('L6').Value = 18

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('L6').Value = 138
$ws.Range('L7').Value = 493

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range('L2').Value = 23
$ws.Range('L3').Value = 32
$ws.Range('L7').Value = 109

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range('L3').Value = 60
$ws.Range('L7').Value = 184

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range('L6').Value = 59
$ws.Range('L7').Value = 174

$ws = $wb.Worksheets.Item('Loop')
$ws.Range('L2').Value = 58
$ws.Range('L7').Value = 316

$ws = $wb.Worksheets.Item('Mount Greenwood')
$ws.Range('L6').Value = 6
$ws.Range('L7').Value = 18

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range('L4').Value = 15
$ws.Range('L6').Value = 29
$ws.Range('L7').Value = 103

$ws = $wb.Worksheets.Item('New City')
$ws.Range('L2').Value = 103
$ws.Range('L7').Value = 290

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('L2').Value = 151
$ws.Range('L3').Value = 203
$ws.Range('L6').Value = 117
$ws.Range('L7').Value = 520

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range('K4').Value = 17
$ws.Range('L6').Value = 30
$ws.Range('K7').Value = 242
$ws.Range('L7').Value = 121

$ws = $wb.Worksheets.Item('River North')
$ws.Range('L6').Value = 110
$ws.Range('L7').Value = 236

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range('L6').Value = 23
$ws.Range('L7').Value = 103

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range('L2').Value = 54
$ws.Range('L7').Value = 199

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range('L2').Value = 133
$ws.Range('L3').Value = 144
$ws.Range('L7').Value = 400

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range('L3').Value = 52
$ws.Range('L7').Value = 146

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('L3').Value = 314
$ws.Range('L6').Value = 164
$ws.Range('L7').Value = 777

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range('L6').Value = 59
$ws.Range('L7').Value = 215

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range('L2').Value = 49
$ws.Range('L7').Value = 150

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range('L2').Value = 42
$ws.Range('L6').Value = 78
$ws.Range('L7').Value = 190

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range('L2').Value = 77
$ws.Range('L7').Value = 203

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range('L2').Value = 55
$ws.Range('L3').Value = 47
$ws.Range('L7').Value = 171

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range('L3').Value = 107
$ws.Range('L7').Value = 259
